# Refactor api response formatting
# Appends a new row (row 51) of API/database log data to each of the four
# worksheets, mirroring the existing rows (time, length/ID/checksum byte
# strings, and their decimal counterparts).

$wb = $excel.ActiveWorkbook

$dateVal = [double]"45837.43663194445"
$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Add-LogRow {
    param($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i)

    $ws.Cells.Item($row, 1).Value = $dateVal
    $ws.Cells.Item($row, 1).NumberFormat = $dateFmt

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e

    $ws.Cells.Item($row, 6).Value = [double]$f
    $ws.Cells.Item($row, 7).Value = [double]$g
    $ws.Cells.Item($row, 8).Value = [double]$h
    $ws.Cells.Item($row, 9).Value = [double]$i
}

# Sheet 1: DE_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
Add-LogRow $ws1 51 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x64" "0x14" 380 "7.598631275147109e+23" 356 14

# Sheet 2: DE_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
Add-LogRow $ws2 51 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x64" "0xe" 380 "5.68432987514711e+23" 356 14

# Sheet 3: DE_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
Add-LogRow $ws3 51 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x7E" "0x7" 130 "5.68631262647114e+23" 126 7

# Sheet 4: DE_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
Add-LogRow $ws4 51 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x7E" "0x3" 130 "9.85046333984776e+23" 126 3
